$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N3").Value = 1.48
$ws.Range("O3").Value = 2.6
$ws.Range("G4").Value = 3.1
$ws.Range("I4").Value = 2.45
$ws.Range("U4").Value = 15
$ws.Range("W4").Value = 34
$ws.Range("J9").Value = 1.05
$ws.Range("L9").Value = 1.29
$ws.Range("J10").Value = 1.04
$ws.Range("L10").Value = 1.22
$ws.Range("G11").Value = 2.82
$ws.Range("H11").Value = 3.05
$ws.Range("I11").Value = 2.45
$ws.Range("P11").Value = 1.5
$ws.Range("Q11").Value = 2.25
$ws.Range("R11").Value = 1.98
$ws.Range("T11").Value = 7
$ws.Range("U11").Value = 13
$ws.Range("V11").Value = 11
$ws.Range("W11").Value = 35
$ws.Range("X11").Value = 30
$ws.Range("AA11").Value = 6
$ws.Range("AB11").Value = 17.5
$ws.Range("AE11").Value = 6.4
$ws.Range("AF11").Value = 10.75
$ws.Range("AG11").Value = 10
$ws.Range("AH11").Value = 26
$ws.Range("AI11").Value = 24
$ws.Range("AJ11").Value = 45
$ws.Range("J13").Value = 1.07
$ws.Range("L13").Value = 1.4
$ws.Range("N13").Value = 2.25
$ws.Range("O13").Value = 1.62
$ws.Range("M16").Value = 3.25
$ws.Range("R16").Value = 1.88
$ws.Range("S16").Value = 1.82
$ws.Range("G19").Value = 2.25
$ws.Range("H19").Value = 3.45
$ws.Range("I19").Value = 2.85
$ws.Range("V19").Value = 9
$ws.Range("AE19").Value = 11
$ws.Range("AF19").Value = 16
$ws.Range("G21").Value = 2.47
$ws.Range("I21").Value = 2.85
$ws.Range("L21").Value = 1.35
$ws.Range("M21").Value = 2.7
$ws.Range("N21").Value = 2.02
$ws.Range("O21").Value = 1.62
$ws.Range("P21").Value = 1.4
$ws.Range("R21").Value = 1.75
$ws.Range("S21").Value = 1.85
$ws.Range("T21").Value = 7.7
$ws.Range("U21").Value = 12
$ws.Range("V21").Value = 9.25
$ws.Range("W21").Value = 27
$ws.Range("X21").Value = 21
$ws.Range("Z21").Value = 8
$ws.Range("AB21").Value = 14
$ws.Range("AC21").Value = 70
$ws.Range("AD21").Value = 600
$ws.Range("AE21").Value = 8
$ws.Range("AG21").Value = 10.5
$ws.Range("AI21").Value = 27
$ws.Range("AJ21").Value = 37
$ws.Range("R24").Value = 1.5
$ws.Range("L25").Value = 1.33
$ws.Range("M25").Value = 3.25
$ws.Range("N25").Value = 2.05
$ws.Range("O25").Value = 1.75
$ws.Range("O26").Value = 1.8
$ws.Range("I27").Value = 3.45
$ws.Range("R27").Value = 2.18
$ws.Range("T27").Value = 5.5
$ws.Range("U27").Value = 9
$ws.Range("AB27").Value = 19.5
$ws.Range("AE27").Value = 7.6
$ws.Range("AF27").Value = 16.5
$ws.Range("H29").Value = 5
$ws.Range("I29").Value = 2.88
$ws.Range("L29").Value = 1.03
$ws.Range("M29").Value = 15
$ws.Range("N29").Value = 1.14
$ws.Range("O29").Value = 5.5
$ws.Range("P29").Value = 1.11
$ws.Range("Q29").Value = 6.5
$ws.Range("R29").Value = 1.18
$ws.Range("S29").Value = 4.5
$ws.Range("T29").Value = 29
$ws.Range("U29").Value = 23
$ws.Range("AA29").Value = 17
$ws.Range("AC29").Value = 19
$ws.Range("G32").Value = 2.5
$ws.Range("I32").Value = 2.7
$ws.Range("K32").Value = 12
$ws.Range("T32").Value = 10
$ws.Range("V32").Value = 10
$ws.Range("W32").Value = 26
$ws.Range("AH32").Value = 26
$ws.Range("G33").Value = 3.25
$ws.Range("I33").Value = 2.2
$ws.Range("R33").Value = 2.1
$ws.Range("S33").Value = 1.67
$ws.Range("V33").Value = 12
$ws.Range("AH33").Value = 21
